# Update '想去人数' (F column) counts per gh-pages regeneration diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 264
$ws.Cells.Item(5, 6).Value = 321
$ws.Cells.Item(6, 6).Value = 459
$ws.Cells.Item(7, 6).Value = 2105
$ws.Cells.Item(9, 6).Value = 46
$ws.Cells.Item(10, 6).Value = 1619
$ws.Cells.Item(11, 6).Value = 1619
$ws.Cells.Item(12, 6).Value = 1351
$ws.Cells.Item(17, 6).Value = 549
$ws.Cells.Item(18, 6).Value = 152
$ws.Cells.Item(19, 6).Value = 3
$ws.Cells.Item(20, 6).Value = 7175
$ws.Cells.Item(21, 6).Value = 7836
$ws.Cells.Item(22, 6).Value = 45
$ws.Cells.Item(24, 6).Value = 190
$ws.Cells.Item(26, 6).Value = 492
$ws.Cells.Item(27, 6).Value = 91
$ws.Cells.Item(31, 6).Value = 16
$ws.Cells.Item(33, 6).Value = 194
$ws.Cells.Item(35, 6).Value = 1421
$ws.Cells.Item(36, 6).Value = 150
$ws.Cells.Item(37, 6).Value = 222
$ws.Cells.Item(39, 6).Value = 289
$ws.Cells.Item(40, 6).Value = 8
$ws.Cells.Item(41, 6).Value = 714
$ws.Cells.Item(43, 6).Value = 1359
$ws.Cells.Item(44, 6).Value = 333
$ws.Cells.Item(45, 6).Value = 240
$ws.Cells.Item(46, 6).Value = 190
$ws.Cells.Item(47, 6).Value = 84
$ws.Cells.Item(48, 6).Value = 159
$ws.Cells.Item(49, 6).Value = 155

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 53
$ws.Cells.Item(18, 6).Value = 296

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 2614
$ws.Cells.Item(4, 6).Value = 279
$ws.Cells.Item(5, 6).Value = 132

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 279
$ws.Cells.Item(5, 6).Value = 132
$ws.Cells.Item(7, 6).Value = 321
$ws.Cells.Item(9, 6).Value = 459
$ws.Cells.Item(10, 6).Value = 2105
$ws.Cells.Item(12, 6).Value = 46
$ws.Cells.Item(13, 6).Value = 1619
$ws.Cells.Item(14, 6).Value = 1619
$ws.Cells.Item(16, 6).Value = 1351
$ws.Cells.Item(20, 6).Value = 549
$ws.Cells.Item(22, 6).Value = 7175
$ws.Cells.Item(23, 6).Value = 7837
$ws.Cells.Item(24, 6).Value = 45
$ws.Cells.Item(26, 6).Value = 190
$ws.Cells.Item(27, 6).Value = 91
$ws.Cells.Item(29, 6).Value = 16
$ws.Cells.Item(30, 6).Value = 194
$ws.Cells.Item(31, 6).Value = 1421
$ws.Cells.Item(32, 6).Value = 150
$ws.Cells.Item(33, 6).Value = 222
$ws.Cells.Item(35, 6).Value = 289
$ws.Cells.Item(38, 6).Value = 714
$ws.Cells.Item(42, 6).Value = 333
$ws.Cells.Item(43, 6).Value = 240
$ws.Cells.Item(44, 6).Value = 190
$ws.Cells.Item(45, 6).Value = 84
$ws.Cells.Item(46, 6).Value = 159
$ws.Cells.Item(47, 6).Value = 155
$ws.Cells.Item(49, 6).Value = 296
